# Testes GJ 1132 b
# Update the exoplanet parameter sheet with the new test-case values and
# refresh the description text for "semiEixoRaioStar".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New input values for the GJ 1132 b test case (row 2 holds the live values
# that correspond to the variable names / headers in row 1).
$ws.Range("E2").Value = 0.21
$ws.Range("M2").Value = 0.015
$ws.Range("N2").Value = 0.18
$ws.Range("O2").Value = 0.101
$ws.Range("P2").Value = 0.0052
$ws.Range("Q2").Value = 1.628
$ws.Range("R2").Value = 86.58
$ws.Range("V2").Value = 1

# Clarify the description of "semiEixoRaioStar" (column A instructions).
$ws.Range("A23").Value = "semiEixoRaioStar = conversão do semi-eixo orbital em relação ao raio da estrela [em UA]"

# Move the active selection to reflect where the user left off editing.
$ws.Range("M21").Select() | Out-Null
